# Trade #36 closed at 2026-02-16 22:55:42 - base_strategy DOWN +0.000%
# Appends a new trade row (row 37) to both the "All Trades" sheet and the
# "base_strategy" sheet, mirroring the existing OPEN-trade row layout.
#
# Row 37 is identical to row 36 in every column except "Trade #" (A) and
# "Time" (C), so the existing row 36 is duplicated via Range.Copy (this
# keeps the blank Exit Price / Exit Reason cells (G, P) truly blank,
# exactly like every other still-OPEN trade row) and then the two changed
# cells are overwritten.

$wb = $excel.ActiveWorkbook

$sheetNames = @("All Trades", "base_strategy")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $srcRow = 36
    $dstRow = 37

    $srcRange = "A" + $srcRow + ":Q" + $srcRow
    $dstRange = "A" + $dstRow + ":Q" + $dstRow

    $ws.Range($srcRange).Copy($ws.Range($dstRange))

    $ws.Cells.Item($dstRow, 1).Value = 36
    $ws.Cells.Item($dstRow, 3).Value = "22:55:42"
}
